$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$passwords = @{
    2  = "ieFSsp3"
    3  = "Kmj9vwJ"
    4  = "1i1V7Pe"
    5  = "rw3yfK1"
    6  = "PWQwDao"
    7  = "be6dvf0"
    8  = "hDnYltE"
    9  = "VqlpGeK"
    10 = "nv9pgB2"
    11 = "1jcEQSN"
    12 = "tejD45z"
    13 = "tc3zGGJ"
    14 = "4H9J2Ld"
    15 = "nzWSjIr"
    16 = "4m4xmnB"
    17 = "4gfCO70"
    18 = "adHaT4C"
    19 = "sHBuJAm"
    20 = "4Vfbka7"
    21 = "6BjEsGI"
    22 = "SRhWOjT"
    23 = "IpBWhcY"
    24 = "69U7n0v"
    25 = "0JiZykn"
    26 = "h0Wmll2"
    27 = "4ejEnqL"
    28 = "EGccYU4"
    29 = "g72C1S7"
    30 = "asfU8rD"
    31 = "qkvO4lr"
}

foreach ($row in $passwords.Keys) {
    $ws.Range("D$row").Value = $passwords[$row]
}
